$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich runs) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "47"

$dateCell = $ws.Range("C9")
$dateCell.Characters(27, 10).Text = "11/20/2023"
$dateCell.Characters(48, 10).Text = "11/26/2023"

# --- Style-changing cells: copy format+value from a same-style reference cell, then set final value ---
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("I14").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 2
$ws.Range("K14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100
$ws.Range("I14").Copy($ws.Range("G23"))
$ws.Range("G23").Value = 2
$ws.Range("K14").Copy($ws.Range("H23"))
$ws.Range("H23").Value = -50
$ws.Range("C14").Copy($ws.Range("F26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("F28"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("I14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("H15").Value = -100
$ws.Range("N15").Value = -38.888888888888
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 160
$ws.Range("J16").Value = 179
$ws.Range("K16").Value = -10.614525139664
$ws.Range("L16").Value = -7.514450867052
$ws.Range("M16").Value = -3.614457831325
$ws.Range("N16").Value = -83.367983367983
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -71.428571428571
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -41.176470588235
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 223
$ws.Range("K17").Value = -8.071748878923
$ws.Range("L17").Value = 3.015075376884
$ws.Range("M17").Value = 61.417322834645
$ws.Range("N17").Value = -23.507462686567
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 7.142857142857
$ws.Range("I18").Value = 200
$ws.Range("J18").Value = 257
$ws.Range("K18").Value = -22.178988326848
$ws.Range("L18").Value = -3.846153846153
$ws.Range("M18").Value = -12.663755458515
$ws.Range("N18").Value = -86.440677966101
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 136.363636363636
$ws.Range("F19").Value = 92
$ws.Range("G19").Value = 78
$ws.Range("H19").Value = 17.948717948717
$ws.Range("I19").Value = 949
$ws.Range("J19").Value = 944
$ws.Range("K19").Value = 0.529661016949
$ws.Range("L19").Value = 30.357142857142
$ws.Range("M19").Value = -28.913857677902
$ws.Range("N19").Value = -62.251392203659
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 57
$ws.Range("J20").Value = 67
$ws.Range("K20").Value = -14.925373134328
$ws.Range("L20").Value = 1.785714285714
$ws.Range("M20").Value = 14
$ws.Range("N20").Value = -95.148936170212
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 17.241379310344
$ws.Range("F21").Value = 132
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = 5.6
$ws.Range("I21").Value = 1583
$ws.Range("J21").Value = 1689
$ws.Range("K21").Value = -6.275902901124
$ws.Range("L21").Value = 14.131218457101
$ws.Range("M21").Value = -17.595002602811
$ws.Range("N21").Value = -75.334995325646
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 83
$ws.Range("K22").Value = -3.614457831325
$ws.Range("M22").Value = 14.285714285714
$ws.Range("F23").Value = 1
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = -35.294117647058
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 75
$ws.Range("F24").Value = 228
$ws.Range("G24").Value = 153
$ws.Range("H24").Value = 49.019607843137
$ws.Range("I24").Value = 2057
$ws.Range("J24").Value = 2126
$ws.Range("K24").Value = -3.245531514581
$ws.Range("L24").Value = 14.150943396226
$ws.Range("M24").Value = 20.645161290322
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 10.810810810810
$ws.Range("I25").Value = 467
$ws.Range("J25").Value = 468
$ws.Range("K25").Value = -0.213675213675
$ws.Range("L25").Value = 12.530120481927
$ws.Range("M25").Value = 21.298701298701
$ws.Range("H26").Value = -100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 4
$ws.Range("H28").Value = -100
$ws.Range("H29").Value = -100
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 400
$ws.Range("I30").Value = 13
$ws.Range("K30").Value = -40.909090909090
$ws.Range("L30").Value = -23.529411764705
